$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1498
$ws.Range("F3").Value = 1465
$ws.Range("F4").Value = 404
$ws.Range("F5").Value = 231
$ws.Range("F6").Value = 730
$ws.Range("F7").Value = 43
$ws.Range("F8").Value = 664
$ws.Range("F11").Value = 1390
$ws.Range("F12").Value = 36539
$ws.Range("G12").Value = "暂时售罄"
$ws.Range("F13").Value = 7245
$ws.Range("F15").Value = 381
$ws.Range("F16").Value = 589
$ws.Range("F17").Value = 457
$ws.Range("F20").Value = 347
$ws.Range("F21").Value = 53
$ws.Range("F22").Value = 460
$ws.Range("F23").Value = 122
$ws.Range("F24").Value = 826
$ws.Range("F25").Value = 22
$ws.Range("F26").Value = 328
$ws.Range("F27").Value = 400
$ws.Range("F28").Value = 452
$ws.Range("F30").Value = 227
$ws.Range("F31").Value = 59
$ws.Range("F32").Value = 752
$ws.Range("F33").Value = 297
$ws.Range("F34").Value = 138
$ws.Range("F35").Value = 764
$ws.Range("F36").Value = 118
$ws.Range("F38").Value = 815
$ws.Range("F39").Value = 297
$ws.Range("F40").Value = 55
$ws.Range("F41").Value = 28
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 1221
$ws.Range("F5").Value = 170
$ws.Range("F6").Value = 295
$ws.Range("F12").Value = 60
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 1485
$ws.Range("F3").Value = 366
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1485
$ws.Range("F3").Value = 366
$ws.Range("F4").Value = 1221
$ws.Range("F5").Value = 1498
$ws.Range("F7").Value = 1465
$ws.Range("F8").Value = 231
$ws.Range("F9").Value = 730
$ws.Range("F10").Value = 43
$ws.Range("F11").Value = 664
$ws.Range("F13").Value = 36539
$ws.Range("F14").Value = 170
$ws.Range("F15").Value = 295
$ws.Range("F20").Value = 7245
$ws.Range("F21").Value = 381
$ws.Range("F22").Value = 60
$ws.Range("F23").Value = 589
$ws.Range("F24").Value = 457
$ws.Range("F27").Value = 347
$ws.Range("F29").Value = 53
$ws.Range("F31").Value = 460
$ws.Range("F32").Value = 122
$ws.Range("F33").Value = 826
$ws.Range("F34").Value = 22
$ws.Range("F35").Value = 328
$ws.Range("F36").Value = 400
$ws.Range("F37").Value = 452
$ws.Range("F39").Value = 227
$ws.Range("F40").Value = 59
$ws.Range("F41").Value = 752
$ws.Range("F43").Value = 297
$ws.Range("F44").Value = 138
$ws.Range("F45").Value = 815
$ws.Range("F46").Value = 297
$ws.Range("F47").Value = 55
$ws.Range("F49").Value = 28
